$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns G (AttackRange) and H (SightRange)
$ws.Range("G1").Value = "AttackRange"
$ws.Range("H1").Value = "SightRange"

# New data values for the four soldier rows
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 500

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 500

$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 500

$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 500

# Column widths: G grows slightly (bestFit), H is a new bestFit column
# (ColumnWidth is in character units; stored OOXML width = ColumnWidth + 5/MaxDigitWidth(7),
#  so subtract 5/7 here to land on the target stored widths of 12.5 and 11.375)
$ws.Columns.Item(7).ColumnWidth = 12.5 - (5/7)
$ws.Columns.Item(8).ColumnWidth = 11.375 - (5/7)

# Update the selection to match the post-edit cursor location
$ws.Range("J8").Select()
